$d = $word.ActiveDocument

# 1. Update the "Last Modified" date field result text.
$d.Content.Find.Execute("2017-03-04", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2017-03-05", 2) | Out-Null

# 2. Add a new bullet point after "... in MainActivity" describing the
#    Score/Tiles layout update, with a manual line break before the TODO note.

# Locate the paragraph that ends with "in MainActivity" (the last bullet
# under the "Implemented" list) so the new bullet is inserted right after it,
# before the trailing blank paragraph / "Notes (REMOVE)" section.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text
    if ($paraText -like "*in MainActivity*") {
        $targetIndex = $i
    }
}

$anchorPara = $d.Paragraphs.Item($targetIndex)

# Inserting a paragraph after the anchor clones its paragraph formatting
# (ListParagraph style, numId 7 bullet list at level 0).
$anchorPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Updated ‘Score’ and ‘Tiles’ activity layouts to use images in header row."

# Insert a manual line break followed by the TODO text right before the
# paragraph's trailing paragraph mark.
$insertPos = $newPara.Range.End - 1
$breakRange = $d.Range($insertPos, $insertPos)
$breakRange.InsertAfter([char]11 + "(TODO: want to transpose row/column for portrait mode but haven’t been able to find a way to do this)")
